# Sanitiza template e remove planilhas de teste
#
# This script scrubs the filled-in answer text from "Planilha Base.xlsx" back
# down to its bare '<n>*...<n>*' placeholder markers, collapses the verbose
# 'SIM, porque ...' / 'ATENDE, pois ...' justifications down to the bare
# 'SIM' / 'ATENDE.' keywords, relabels a few header cells in row 13, and
# (re)populates the second 'META ESPECÍFICA' template block in row 25 with the
# same sanitized placeholder content as row 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch A1 so the sheet's used range (dimension) grows to include row 1,
# mirroring the source edit (A2:L137 -> A1:L137) without altering its
# appearance (re-applying the default 'Normal' style is a no-op visually).
$ws.Range("A1").Style = "Normal"

# DIAGNÓSTICO row (row 4): drop the filled-in justification, keep the cell
# (and its style) but make it blank.
$ws.Range("F4").Value = ""

# META GERAL row (row 8): collapse the long narrative text in the title
# cell down to its bare placeholder markers, and all the justification
# cells down to the single word 'SIM'.
$ws.Range("A8").Value = "1*1*"
$ws.Range("B8").Value = "SIM"
$ws.Range("C8").Value = "SIM"
$ws.Range("D8").Value = "SIM"
$ws.Range("E8").Value = "SIM"
$ws.Range("F8").Value = "SIM"

# Fórmula / Referência utilizada (row 10): same placeholder-marker collapse.
$ws.Range("F10").Value = "0*0*"

# ESTRATÉGIA DE IMPLEMENTAÇÃO row (row 12): collapse justifications to
# the bare 'ATENDE.' keyword.
$ws.Range("B12").Value = "ATENDE."
$ws.Range("C12").Value = "ATENDE."
$ws.Range("D12").Value = "ATENDE."
$ws.Range("E12").Value = "ATENDE."

# META ESPECÍFICA header row (row 13): relabel the PESP/PNSP/MJSP-specific
# column headers to their generic 'setorial/nacional/informada' wording.
$ws.Range("G13").Value = "A Meta setorial foi informada? Existe aderência?"
$ws.Range("H13").Value = "A Meta nacional foi informada? Existe aderência?"
$ws.Range("I13").Value = "A política informada foi apresentada? Existe aderência?"

# META ESPECÍFICA data row (row 14): collapse the title cell and all the
# embedded '<n>*...<n>*' free-text spans down to bare markers, leaving the
# surrounding boilerplate sentences untouched.
$ws.Range("A14").Value = "2*2*"
$ws.Range("E14").Value = "SIM.`n`nA referência informada foi:`n`n`n`n3*3*"
$ws.Range("F14").Value = "SIM.`n`nO Indicador e Fórmula de Cálculo informado foi:`n`n`n`n4*Descrição do Indicador:`n4*`n`n`n`n5*Fórmula:`n5*`n`nO indicador e a fórmula de cálculo são adequados para o eficiente monitoramento da meta."
$ws.Range("G14").Value = "SIM.`n`nA Meta informada foi:`n`n`n`n6*6*`n`n`n`nExiste aderência da referida Meta à Política Pública."
$ws.Range("H14").Value = "SIM.`n`nA Meta informada foi:`n`n`n7*7*`n`n`n`nExiste aderência da referida Meta Específica à Política informada."
$ws.Range("I14").Value = "SIM.`n`nA política informada foi:`n`n`n`n8*8*`n`nExiste aderência da referida Meta Específica à Política informada."

# Second META ESPECÍFICA template block (row 25) was blank in the original
# workbook (merged placeholder cells with no content) - populate it with the
# same sanitized placeholder content as row 14 so the template offers two
# ready-to-fill blocks.
$ws.Range("A25").Value = "2*2*"
$ws.Range("B25").Value = "SIM"
$ws.Range("C25").Value = "SIM"
$ws.Range("D25").Value = "SIM"
$ws.Range("E25").Value = "SIM.`n`nA referência informada foi:`n`n`n`n3*3*"
$ws.Range("F25").Value = "SIM.`n`nO Indicador e Fórmula de Cálculo informado foi:`n`n`n`n4*Descrição do Indicador:`n4*`n`n`n`n5*Fórmula:`n5*`n`nO indicador e a fórmula de cálculo são adequados para o eficiente monitoramento da meta."
$ws.Range("G25").Value = "SIM.`n`nA Meta informada foi:`n`n`n`n6*6*`n`n`n`nExiste aderência da referida Meta à Política Pública."
$ws.Range("H25").Value = "SIM.`n`nA Meta informada foi:`n`n`n7*7*`n`n`n`nExiste aderência da referida Meta Específica à Política informada."
$ws.Range("I25").Value = "SIM.`n`nA política informada foi:`n`n`n`n8*8*`n`nExiste aderência da referida Meta Específica à Política informada."

